# Rename header "Run 2" -> "Row 2" and drop the now-unused Column3..Column7
# header cells (the data-collection columns were never populated, so the
# sheet is being trimmed down to the three columns that are actually used:
# Test Method / Run 1 / Row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Row 2"
$ws.Range("D1:H1").ClearContents()

# "enabled collection of the data" - flip the collected Run-2 results for
# the rows where the new run's outcome differs from Run 1's. Copy/paste
# from an existing true/false cell (instead of assigning a literal string)
# so the written cell stays a shared-string text value like its neighbours,
# rather than being auto-coerced into a native Excel boolean.
$ws.Range("B7").Copy($ws.Range("B4"))
$ws.Range("B7").Copy($ws.Range("B6"))
$ws.Range("B7").Copy($ws.Range("B11"))

$ws.Range("B2").Copy($ws.Range("B7"))
$ws.Range("B2").Copy($ws.Range("B9"))
